# Refresh the crypto price / 1h-volume columns (D, E) for rows 2-51
# with the latest scrape values. Each D-column price is prefixed with a
# leading apostrophe so Excel stores it as literal text (quotePrefix)
# instead of re-parsing look-alike numbers such as "25.920.04" or
# "4.271" into a Double and silently dropping the exact formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$quote = "'"

$ws.Range("D2").Value = $quote + '25.920.04'
$ws.Range("E2").Value = '  +0.06%  '
$ws.Range("D3").Value = $quote + '1.633.66'
$ws.Range("E3").Value = '  -0.08%  '
$ws.Range("D4").Value = $quote + '1.002'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = $quote + '216.14'
$ws.Range("E5").Value = '  +0.79%  '
$ws.Range("D6").Value = $quote + '0.5099'
$ws.Range("E6").Value = '  +0.20%  '
$ws.Range("D7").Value = $quote + '1.002'
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = $quote + '0.2580'
$ws.Range("E8").Value = '  +1.14%  '
$ws.Range("D9").Value = $quote + '0.06358'
$ws.Range("E9").Value = '  +0.25%  '
$ws.Range("D10").Value = $quote + '19.48'
$ws.Range("E10").Value = '  +0.14%  '
$ws.Range("D11").Value = $quote + '0.07776'
$ws.Range("E11").Value = '  +0.17%  '
$ws.Range("D12").Value = $quote + '4.271'
$ws.Range("E12").Value = '  -0.14%  '
$ws.Range("D13").Value = $quote + '1.636.43'
$ws.Range("E13").Value = '  -0.07%  '
$ws.Range("D14").Value = $quote + '1.859.19'
$ws.Range("E14").Value = '  -0.06%  '
$ws.Range("D15").Value = $quote + '0.5514'
$ws.Range("E15").Value = '  +1.71%  '
$ws.Range("D16").Value = $quote + '63.99'
$ws.Range("E16").Value = '  -0.13%  '
$ws.Range("D17").Value = $quote + '0.0₅7673'
$ws.Range("E17").Value = '  -0.38%  '
$ws.Range("D18").Value = $quote + '25.946.09'
$ws.Range("E18").Value = '  +0.06%  '
$ws.Range("D19").Value = $quote + '1.002'
$ws.Range("E19").Value = '  -0.07%  '
$ws.Range("D20").Value = $quote + '196.09'
$ws.Range("E20").Value = '  +0.35%  '
$ws.Range("D21").Value = $quote + '4.417'
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("D22").Value = $quote + '9.909'
$ws.Range("E22").Value = '  +0.17%  '
$ws.Range("D23").Value = $quote + '6.065'
$ws.Range("E23").Value = '  +0.88%  '
$ws.Range("D24").Value = $quote + '1.003'
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("D25").Value = $quote + '1.906'
$ws.Range("E25").Value = '  +2.19%  '
$ws.Range("D26").Value = $quote + '142.36'
$ws.Range("E26").Value = '  +1.08%  '
$ws.Range("D27").Value = $quote + '0.1252'
$ws.Range("E27").Value = '  +4.82%  '
$ws.Range("D28").Value = $quote + '15.65'
$ws.Range("E28").Value = '  +0.42%  '
$ws.Range("D29").Value = $quote + '6.768'
$ws.Range("E29").Value = '  -0.61%  '
$ws.Range("D30").Value = $quote + '1.241'
$ws.Range("E30").Value = '  +0.60%  '
$ws.Range("D31").Value = $quote + '0.04915'
$ws.Range("E31").Value = '  +0.35%  '
$ws.Range("D32").Value = $quote + '3.250'
$ws.Range("E32").Value = '  +0.34%  '
$ws.Range("D33").Value = $quote + '3.203'
$ws.Range("E33").Value = '  +1.19%  '
$ws.Range("D34").Value = $quote + '1.543'
$ws.Range("E34").Value = '  +1.25%  '
$ws.Range("D35").Value = $quote + '2.370'
$ws.Range("E35").Value = '  +0.12%  '
$ws.Range("D36").Value = $quote + '0.8985'
$ws.Range("E36").Value = '  +1.08%  '
$ws.Range("D37").Value = $quote + '0.5541'
$ws.Range("E37").Value = '  +2.71%  '
$ws.Range("D38").Value = $quote + '2.539'
$ws.Range("E38").Value = '  -1.65%  '
$ws.Range("D39").Value = $quote + '1.113.10'
$ws.Range("E39").Value = '  -2.32%  '
$ws.Range("D40").Value = $quote + '0.01559'
$ws.Range("E40").Value = '  +0.76%  '
$ws.Range("D41").Value = $quote + '1.001'
$ws.Range("E41").Value = '  -0.12%  '
$ws.Range("D42").Value = $quote + '5.620'
$ws.Range("E42").Value = '  +3.50%  '
$ws.Range("D43").Value = $quote + '0.7954'
$ws.Range("E43").Value = '  -2.17%  '
$ws.Range("D44").Value = $quote + '97.51'
$ws.Range("E44").Value = '  -1.19%  '
$ws.Range("D45").Value = $quote + '0.0₈119'
$ws.Range("E45").Value = '  -6.44%  '
$ws.Range("D46").Value = $quote + '1.770.14'
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("D47").Value = $quote + '0.4447'
$ws.Range("E47").Value = '  -1.78%  '
$ws.Range("D48").Value = $quote + '1.004'
$ws.Range("E48").Value = '  +0.62%  '
$ws.Range("D49").Value = $quote + '54.81'
$ws.Range("E49").Value = '  +0.31%  '
$ws.Range("D50").Value = $quote + '0.05135'
$ws.Range("E50").Value = '  +1.63%  '
$ws.Range("D51").Value = $quote + '7.570'
$ws.Range("E51").Value = '  +3.14%  '
